$d = $word.ActiveDocument

$replacements = @(
    @("2025-09-03 Wednesday", "2025-09-04 Thursday"),
    @("816÷8=", "691÷9="),
    @("778÷9=", "629÷5="),
    @("935÷8=", "552÷8="),
    @("219÷6=", "979÷3="),
    @("983÷6=", "870÷3="),
    @("563÷2=", "947÷2="),
    @("943÷2=", "167÷3="),
    @("554÷5=", "842÷2="),
    @("356÷6=", "674÷5="),
    @("981÷4=", "961÷6="),
    @("515÷8=", "885÷4="),
    @("159÷4=", "402÷5="),
    @("420÷9=", "253÷4="),
    @("795÷7=", "670÷7="),
    @("640÷5=", "359÷9="),
    @("367÷2=", "737÷9="),
    @("561÷8=", "691÷5="),
    @("554÷8=", "297÷5="),
    @("602÷6=", "554÷9="),
    @("544÷9=", "683÷3="),
    @("250÷5=", "604÷9="),
    @("546÷5=", "735÷7="),
    @("729÷4=", "826÷6="),
    @("744÷8=", "273÷6="),
    @("841÷6=", "592÷3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Done applying replacements"
